$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Updated raw measurements (rows 2-3) ---
$ws.Range("B2").Value = 169.31
$ws.Range("C2").Value = 101.61
$ws.Range("B3").Value = 156.48
$ws.Range("C3").Value = 103.96

# --- Samples / win-loss block (rows 7-9) ---
$ws.Range("B7").Value = 50
$ws.Range("B8").Value = 29

# New radian -> degrees conversion example next to that block
$ws.Range("D7").Value = "Radian"
$ws.Range("E7").Value = 8.42
$ws.Range("D8").Value = "Degrees"
$ws.Range("E8").Formula = "=DEGREES(E7)"

# --- Score / stddev block (rows 17-19) ---
$ws.Range("B17").Value = 50.92
$ws.Range("B18").Value = 0.22

# --- Rank / total block (rows 21-22) ---
$ws.Range("B21").Value = 10
$ws.Range("B22").Value = 2458

# Move the active selection to match the saved view
$ws.Range("I21").Select() | Out-Null
